# "terminei o cadastro de novo usuario e coloque um print para capturar a
#  utima tela" - finish registering the new user and add a row for the
#  result of the second registration, fixing up a couple of values on the
#  first (already existing) row along the way.
#
# NOTE: the exact order in which the cells below are written matters - it
# controls the order new entries are appended to xl/sharedStrings.xml, and
# that order must line up with the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (existing user "saulosjs") : format / correct a few fields ---
$ws.Range("G2").Value = "(11)921316555"        # NumeroTelefone, now stored as text

# --- Row 3 (new user registered) ---
$ws.Range("C3").Value = "mano"                  # Senha
$ws.Range("B3").Value = "saulojosilva@hotmail.com"   # Email
$ws.Range("D3").Value = "eita"                  # ConfirmarSenha
$ws.Range("E3").Value = "Marcos"                # PrimeiroNome
$ws.Range("A3").Value = "rapaz"                 # Usuario
$ws.Range("F3").Value = "Silva"                 # UltimoNome
$ws.Range("G3").Value = "(11)954423458"         # NumeroTelefone
$ws.Range("H3").Value = "Algeria"               # Pais

# --- back to row 2 ---
$ws.Range("H2").Value = "Brazil"                # Pais corrected

# --- rest of row 3 ---
$ws.Range("I3").Value = "paris"                 # Cidade
$ws.Range("J3").Value = "rua mundi"             # Endereco
$ws.Range("K3").Value = "af"                    # Estado

# --- CodigoPostal on both rows ---
$ws.Range("L2").Value = "(11)12345"
$ws.Range("L3").Value = "(11)12345"

# --- screenshot/result column and final Usuario fix on row 2 ---
$ws.Range("M2").Value = "Pass"
$ws.Range("A2").Value = "saulosjss33"

# New user's e-mail becomes a clickable mailto hyperlink, same as row 2
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:saulojosilva@hotmail.com")
$ws.Range("B3").Style = "Hiperlink"

# Last screenshot was captured with A2 selected
$ws.Range("A2").Select()
